$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.866.94'
$ws.Range('E2').Value = '  -2.33%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.753.01'
$ws.Range('E3').Value = '  -4.63%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.96'
$ws.Range('E5').Value = '  -8.50%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9989'
$ws.Range('E6').Value = '  -0.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5096'
$ws.Range('E7').Value = '  -5.36%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '42.34'
$ws.Range('E8').Value = '  -5.60%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2789'
$ws.Range('E9').Value = '  -6.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06207'
$ws.Range('E10').Value = '  -10.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.745.34'
$ws.Range('E11').Value = '  -5.07%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.82'
$ws.Range('E12').Value = '  -9.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.06964'
$ws.Range('E13').Value = '  -3.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6113'
$ws.Range('E14').Value = '  -16.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.526'
$ws.Range('E15').Value = '  -9.32%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '77.67'
$ws.Range('E16').Value = '  -12.79%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.000'
$ws.Range('E17').Value = '  -0.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9995'
$ws.Range('E18').Value = '  -0.09%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '25.872.74'
$ws.Range('E19').Value = '  -2.37%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000006967'
$ws.Range('E20').Value = '  -11.96%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.69'
$ws.Range('E21').Value = '  -15.35%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.969.47'
$ws.Range('E22').Value = '  -5.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.087'
$ws.Range('E23').Value = '  -10.85%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.269'
$ws.Range('E24').Value = '  -12.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.230'
$ws.Range('E25').Value = '  -10.56%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '137.69'
$ws.Range('E26').Value = '  -3.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.482'
$ws.Range('E27').Value = '  -13.35%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.826'
$ws.Range('E28').Value = '  -15.84%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.04'
$ws.Range('E29').Value = '  -11.27%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '103.79'
$ws.Range('E30').Value = '  -6.47%  '
$ws.Range('E31').Value = '  -7.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.712'
$ws.Range('E32').Value = '  -12.67%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.516'
$ws.Range('E33').Value = '  -12.88%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04529'
$ws.Range('E34').Value = '  -6.59%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9982'
$ws.Range('E35').Value = '  -0.13%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.634'
$ws.Range('E36').Value = '  -9.75%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9928'
$ws.Range('E37').Value = '  -12.13%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.6145'
$ws.Range('E38').Value = '  -15.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.696'
$ws.Range('E39').Value = '  -12.85%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01557'
$ws.Range('E40').Value = '  -8.86%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9989'
$ws.Range('E41').Value = '  -0.13%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.903'
$ws.Range('E42').Value = '  -17.14%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '103.64'
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3883'
$ws.Range('E44').Value = '  -17.27%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.7468'
$ws.Range('E45').Value = '  -17.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.933'
$ws.Range('E46').Value = '  -15.96%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05410'
$ws.Range('E47').Value = '  -6.08%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1115'
$ws.Range('E48').Value = '  -10.56%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.013'
$ws.Range('E49').Value = '  -18.73%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '30.24'
$ws.Range('E50').Value = '  -12.95%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '52.88'
$ws.Range('E51').Value = '  -12.03%  '
